$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 149.656361
$ws.Range("H2").Value = 448.969083
$ws.Range("I2").Value = 0.5921360794347563
$ws.Range("J2").Value = 0.5921360794347564
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 147.4213356666667
$ws.Range("N2").Value = 442.264007
$ws.Range("O2").Value = 0.9507885170992249
$ws.Range("P2").Value = 0.950788517099225
$ws.Range("Q2").Value = 22062.54062963284
$ws.Range("R2").Value = 198562.8656666956
$ws.Range("S2").Value = 0.5629961848867208
$ws.Range("T2").Value = 0.562996184886721

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 149.656361
$ws.Range("H3").Value = 448.969083
$ws.Range("I3").Value = 0.5921360794347563
$ws.Range("J3").Value = 0.5921360794347564
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 2.340788333333334
$ws.Range("N3").Value = 7.022365000000001
$ws.Range("O3").Value = 0.01509682881537204
$ws.Range("P3").Value = 0.01509682881537204
$ws.Range("Q3").Value = 350.3138638379216
$ws.Range("R3").Value = 3152.824774541295
$ws.Range("S3").Value = 0.008939377026632054
$ws.Range("T3").Value = 0.008939377026632056

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 149.656361
$ws.Range("H4").Value = 448.969083
$ws.Range("I4").Value = 0.5921360794347563
$ws.Range("J4").Value = 0.5921360794347564
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 5.289533666666667
$ws.Range("N4").Value = 15.868601
$ws.Range("O4").Value = 0.03411465408540306
$ws.Range("P4").Value = 0.03411465408540307
$ws.Range("Q4").Value = 791.6123599403203
$ws.Range("R4").Value = 7124.511239462883
$ws.Range("S4").Value = 0.02020051752140346
$ws.Range("T4").Value = 0.02020051752140347

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 52.73412466666667
$ws.Range("H5").Value = 158.202374
$ws.Range("I5").Value = 0.208649853730866
$ws.Range("J5").Value = 0.208649853730866
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 147.4213356666667
$ws.Range("N5").Value = 442.264007
$ws.Range("O5").Value = 0.9507885170992249
$ws.Range("P5").Value = 0.950788517099225
$ws.Range("Q5").Value = 7774.135093572513
$ws.Range("R5").Value = 69967.21584215261
$ws.Range("S5").Value = 0.1983818850217403
$ws.Range("T5").Value = 0.1983818850217403

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 52.73412466666667
$ws.Range("H6").Value = 158.202374
$ws.Range("I6").Value = 0.208649853730866
$ws.Range("J6").Value = 0.208649853730866
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 2.340788333333334
$ws.Range("N6").Value = 7.022365000000001
$ws.Range("O6").Value = 0.01509682881537204
$ws.Range("P6").Value = 0.01509682881537204
$ws.Range("Q6").Value = 123.4394237882789
$ws.Range("R6").Value = 1110.95481409451
$ws.Range("S6").Value = 0.003149951124127298
$ws.Range("T6").Value = 0.003149951124127299

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 52.73412466666667
$ws.Range("H7").Value = 158.202374
$ws.Range("I7").Value = 0.208649853730866
$ws.Range("J7").Value = 0.208649853730866
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 5.289533666666667
$ws.Range("N7").Value = 15.868601
$ws.Range("O7").Value = 0.03411465408540306
$ws.Range("P7").Value = 0.03411465408540307
$ws.Range("Q7").Value = 278.9389278065305
$ws.Range("R7").Value = 2510.450350258774
$ws.Range("S7").Value = 0.007118017584998439
$ws.Range("T7").Value = 0.007118017584998441

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 50.34932566666667
$ws.Range("H8").Value = 151.047977
$ws.Range("I8").Value = 0.1992140668343777
$ws.Range("J8").Value = 0.1992140668343777
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 147.4213356666667
$ws.Range("N8").Value = 442.264007
$ws.Range("O8").Value = 0.9507885170992249
$ws.Range("P8").Value = 0.950788517099225
$ws.Range("Q8").Value = 7422.564839695981
$ws.Range("R8").Value = 66803.08355726385
$ws.Range("S8").Value = 0.1894104471907638
$ws.Range("T8").Value = 0.1894104471907639

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 50.34932566666667
$ws.Range("H9").Value = 151.047977
$ws.Range("I9").Value = 0.1992140668343777
$ws.Range("J9").Value = 0.1992140668343777
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 2.340788333333334
$ws.Range("N9").Value = 7.022365000000001
$ws.Range("O9").Value = 0.01509682881537204
$ws.Range("P9").Value = 0.01509682881537204
$ws.Range("Q9").Value = 117.8571141117339
$ws.Range("R9").Value = 1060.714027005605
$ws.Range("S9").Value = 0.003007500664612683
$ws.Range("T9").Value = 0.003007500664612684

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 50.34932566666667
$ws.Range("H10").Value = 151.047977
$ws.Range("I10").Value = 0.1992140668343777
$ws.Range("J10").Value = 0.1992140668343777
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 5.289533666666667
$ws.Range("N10").Value = 15.868601
$ws.Range("O10").Value = 0.03411465408540306
$ws.Range("P10").Value = 0.03411465408540307
$ws.Range("Q10").Value = 266.3244532077975
$ws.Range("R10").Value = 2396.920078870177
$ws.Range("S10").Value = 0.006796118979001161
$ws.Range("T10").Value = 0.006796118979001163

Write-Output "Applied Thbs1-Itga6 updates (Dr Hou advice)"